$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Range("A2").Value = "Pipeline(steps=[('scaler', None),`n                ('selector',`n                 <__main__.NamedFeatureSelector object at 0x7f91043e86a0>),`n                ('model',`n                 XGBClassifier(base_score=None, booster=None, callbacks=None,`n                               colsample_bylevel=None, colsample_bynode=None,`n                               colsample_bytree=0.8, early_stopping_rounds=None,`n                               enable_categorical=False, eval_metric=None,`n                               feature_types=None, gamma=0.1, gpu_id=None,`n                               grow_policy=None, importance_type=None,`n                               interaction_constraints=None, learning_rate=0.01,`n                               max_bin=None, max_cat_threshold=None,`n                               max_cat_to_onehot=None, max_delta_step=None,`n                               max_depth=3, max_leaves=None,`n                               min_child_weight=None, missing=nan,`n                               monotone_constraints=None, n_estimators=200,`n                               n_jobs=None, num_parallel_tree=None,`n                               predictor=None, random_state=42, ...))])"
$ws.Range("C2").Value = "{'selector': <__main__.NamedFeatureSelector object at 0x7f91043e81c0>, 'scaler': None, 'model__subsample': 0.8, 'model__n_estimators': 200, 'model__max_depth': 3, 'model__learning_rate': 0.01, 'model__gamma': 0.1, 'model__colsample_bytree': 0.8}"
$ws.Range("H2").Value = 0.9777132651805929
$ws.Range("I2").Value = 0.005753295364970079
$ws.Range("J2").Value = 0.523993860387978
$ws.Range("K2").Value = 0.1612873099309589

# --- Row 3 ---
$ws.Range("A3").Value = "Pipeline(steps=[('scaler', RobustScaler()),`n                ('selector',`n                 <__main__.NamedFeatureSelector object at 0x7f91043e87c0>),`n                ('model',`n                 XGBClassifier(base_score=None, booster=None, callbacks=None,`n                               colsample_bylevel=None, colsample_bynode=None,`n                               colsample_bytree=0.5, early_stopping_rounds=None,`n                               enable_categorical=False, eval_metric=None,`n                               feature_types=None, gamma=0, gpu_id=None,`n                               grow_policy=None, importance_type=None,`n                               interaction_constraints=None, learning_rate=0.01,`n                               max_bin=None, max_cat_threshold=None,`n                               max_cat_to_onehot=None, max_delta_step=None,`n                               max_depth=5, max_leaves=None,`n                               min_child_weight=None, missing=nan,`n                               monotone_constraints=None, n_estimators=200,`n                               n_jobs=None, num_parallel_tree=None,`n                               predictor=None, random_state=42, ...))])"
$ws.Range("C3").Value = "{'selector': <__main__.NamedFeatureSelector object at 0x7f91043e8c10>, 'scaler': RobustScaler(), 'model__subsample': 0.8, 'model__n_estimators': 200, 'model__max_depth': 5, 'model__learning_rate': 0.01, 'model__gamma': 0, 'model__colsample_bytree': 0.5}"
$ws.Range("H3").Value = 0.9668544652549609
$ws.Range("I3").Value = 0.006238963318035465
$ws.Range("J3").Value = 0.6465173522555875
$ws.Range("K3").Value = 0.1006492796668859

# --- Row 4 ---
$ws.Range("A4").Value = "Pipeline(steps=[('scaler', RobustScaler()),`n                ('selector',`n                 <__main__.NamedFeatureSelector object at 0x7f91043e8b20>),`n                ('model',`n                 XGBClassifier(base_score=None, booster=None, callbacks=None,`n                               colsample_bylevel=None, colsample_bynode=None,`n                               colsample_bytree=0.5, early_stopping_rounds=None,`n                               enable_categorical=False, eval_metric=None,`n                               feature_types=None, gamma=0, gpu_id=None,`n                               grow_policy=None, importance_type=None,`n                               interaction_constraints=None, learning_rate=0.01,`n                               max_bin=None, max_cat_threshold=None,`n                               max_cat_to_onehot=None, max_delta_step=None,`n                               max_depth=7, max_leaves=None,`n                               min_child_weight=None, missing=nan,`n                               monotone_constraints=None, n_estimators=100,`n                               n_jobs=None, num_parallel_tree=None,`n                               predictor=None, random_state=42, ...))])"
$ws.Range("B4").Value = 0.7347549019607842
$ws.Range("C4").Value = "{'selector': <__main__.NamedFeatureSelector object at 0x7f91043e8e80>, 'scaler': RobustScaler(), 'model__subsample': 0.5, 'model__n_estimators': 100, 'model__max_depth': 7, 'model__learning_rate': 0.01, 'model__gamma': 0, 'model__colsample_bytree': 0.5}"
$ws.Range("H4").Value = 0.9811858061171492
$ws.Range("I4").Value = 0.004929119009328748
$ws.Range("J4").Value = 0.6213583269671505
$ws.Range("K4").Value = 0.1187754442379667

# --- Row 5 ---
$ws.Range("A5").Value = "Pipeline(steps=[('scaler', RobustScaler()),`n                ('selector',`n                 <__main__.NamedFeatureSelector object at 0x7f91043e8eb0>),`n                ('model',`n                 XGBClassifier(base_score=None, booster=None, callbacks=None,`n                               colsample_bylevel=None, colsample_bynode=None,`n                               colsample_bytree=0.5, early_stopping_rounds=None,`n                               enable_categorical=False, eval_metric=None,`n                               feature_types=None, gamma=0, gpu_id=None,`n                               grow_policy=None, importance_type=None,`n                               interaction_constraints=None, learning_rate=0.01,`n                               max_bin=None, max_cat_threshold=None,`n                               max_cat_to_onehot=None, max_delta_step=None,`n                               max_depth=3, max_leaves=None,`n                               min_child_weight=None, missing=nan,`n                               monotone_constraints=None, n_estimators=50,`n                               n_jobs=None, num_parallel_tree=None,`n                               predictor=None, random_state=42, ...))])"
$ws.Range("C5").Value = "{'selector': <__main__.NamedFeatureSelector object at 0x7f91044a5cd0>, 'scaler': RobustScaler(), 'model__subsample': 0.5, 'model__n_estimators': 50, 'model__max_depth': 3, 'model__learning_rate': 0.01, 'model__gamma': 0, 'model__colsample_bytree': 0.5}"
$ws.Range("H5").Value = 0.9723324369519848
$ws.Range("I5").Value = 0.004316199123202302
$ws.Range("J5").Value = 0.6600916204710323
$ws.Range("K5").Value = 0.1241107210684721
